# Atualização de bases das ligas, do dia: 03-06-2024 às 23:01
# Swap the data (columns B through AD) between pairs of rows while keeping
# the sequential index in column A fixed in place. Four pairs of rows were
# reordered in the source data: (73,74), (112,113), (124,125), (159,160).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(73, 74),
    @(112, 113),
    @(124, 125),
    @(159, 160)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
